# Generate Report for Handoff
# Replace the old GUID-based file name (32c1aa73-65ad-4fc1-ab06-b4eb8e77e00b)
# with the new one (d37cc4bd-9942-4a11-af6a-3d42316e1488) everywhere it is
# used, and bump the associated handoff timestamps, across the three
# worksheets: "Overview", "zh-cn", "de-de".

$wb = $excel.ActiveWorkbook

$oldGuid = "32c1aa73-65ad-4fc1-ab06-b4eb8e77e00b"
$newGuid = "d37cc4bd-9942-4a11-af6a-3d42316e1488"

$oldZhHash = "6325feb7507b78b675c038e7bca5c7a4c391f36d"
$newZhHash = "355c77e35c867c8cde17882a04bd037646a72422"

# ---------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-01 13:05:28"

$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# ---------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 13:05:23"

$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"

# ---------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 13:05:28"

$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
